# Auto-generated edit script: updates cryptos list values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.482.03"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").Value = "3.671.54"
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'617.79"
$ws.Range("E5").Value = "  -8.23%  "
$ws.Range("D6").Value = "'159.40"
$ws.Range("E6").Value = "  -1.61%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -0.56%  "
$ws.Range("E9").Value = "  -1.79%  "
$ws.Range("D10").Value = "'7.18"
$ws.Range("E10").Value = "  +1.00%  "
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("E12").Value = "  -2.92%  "
$ws.Range("D13").Value = "4.292.79"
$ws.Range("E13").Value = "  -1.02%  "
$ws.Range("D14").Value = "'32.41"
$ws.Range("E14").Value = "  -1.49%  "
$ws.Range("D15").Value = "3.677.18"
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("D16").Value = "69.540.33"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("D19").Value = "'15.87"
$ws.Range("E19").Value = "  -2.59%  "
$ws.Range("D20").Value = "'10.28"
$ws.Range("E20").Value = "  +4.64%  "
$ws.Range("D21").Value = "'468.73"
$ws.Range("E21").Value = "  -0.97%  "
$ws.Range("D22").Value = "'0.648"
$ws.Range("E22").Value = "  -0.96%  "
$ws.Range("D23").Value = "'79.34"
$ws.Range("E23").Value = "  -1.40%  "
$ws.Range("D24").Value = "3.820.07"
$ws.Range("E24").Value = "  -1.01%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").Value = "'0.0000122"
$ws.Range("E26").Value = "  -4.64%  "
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("D28").Value = "'8.70"
$ws.Range("E28").Value = "  -5.01%  "
$ws.Range("E29").Value = "  -3.45%  "
$ws.Range("D30").Value = "'1.66"
$ws.Range("E30").Value = "  -4.56%  "
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("E32").Value = "  -2.35%  "
$ws.Range("D33").Value = "'26.59"
$ws.Range("E33").Value = "  -1.32%  "
$ws.Range("E34").Value = "  -3.13%  "
$ws.Range("E35").Value = "  -3.31%  "
$ws.Range("D36").Value = "3.673.36"
$ws.Range("E36").Value = "  -0.69%  "
$ws.Range("D37").Value = "'8.27"
$ws.Range("E37").Value = "  -3.41%  "
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").Value = "'178.56"
$ws.Range("E39").Value = "  +2.52%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "'5.77"
$ws.Range("E41").Value = "  -5.66%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'2.21"
$ws.Range("E42").Value = "  -2.21%  "
$ws.Range("E43").Value = "  -2.62%  "
$ws.Range("E44").Value = "  -1.69%  "
$ws.Range("D45").Value = "'46.78"
$ws.Range("E45").Value = "  -0.75%  "
$ws.Range("D46").Value = "'29.19"
$ws.Range("E46").Value = "  +5.14%  "
$ws.Range("D47").Value = "'2.70"
$ws.Range("E47").Value = "  -2.88%  "
$ws.Range("E48").Value = "  -0.39%  "
$ws.Range("D49").Value = "'0.000263"
$ws.Range("E49").Value = "  -6.75%  "
$ws.Range("E50").Value = "  -5.02%  "
$ws.Range("E51").Value = "  -6.81%  "

# Clear the 'quote-prefix' style flag these text-look-alike numbers picked
# up, so the cells end up with the same (default/general) style as before.
# (Applied per-cell: a comma-unioned Range().Style setter only touches the
# first area in this runtime.)
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"

